$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r=12; a=9; b=0.00624999999999999948; c=0; d="Map datatype(Version 4)"; e=20; f=6 },
    @{ r=13; a=10; b=0.00555555555555555577; c=0; d="Map datatype(Version 4)"; e=20; f=6 },
    @{ r=14; a=11; b=0.04930555555555555386; c=0.00069444444444444447; d="Map datatype(Version 4)"; e=20; f=7 },
    @{ r=15; a=12; b=0.00416666666666666661; c=0; d="Map datatype(Version 5)"; e=20; f=6 },
    @{ r=16; a=13; b=0.03125000000000000000; c=0; d="Map datatype(Version 5)"; e=20; f=7 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 2).NumberFormat = "h:mm"
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 3).NumberFormat = "h:mm"
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
}

$ws.Range("C19").Select() | Out-Null
